# Reran driver analyses for nutrients: add 8 new driver rows (alkalinity, EC,
# NO3.N, NO2.N, mineral.N, Tot.N, PO4.P, Tot.P) below the existing PC_axis1 row,
# wire up the same formula pattern used by the existing rows, and restyle the
# blank separator row (row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new driver data (rows 8-15) -----------------------------------------
$newRows = @(
    @{ Row=8;  Name="alkalinity"; B=0.0364520214788158;    C=4.297974;   D="none";        J="mmol/L" },
    @{ Row=9;  Name="EC";         B=0.000913068677324301;  C=541.5156;   D="log10 + 0.1"; J="µS/cm" },
    @{ Row=10; Name="NO3.N";      B=0.0186814009024428;    C=1.879504;   D="log10 + 0.1"; J="mg/L" },
    @{ Row=11; Name="NO2.N";      B=-0.000270790995308241; C=0.01610435; D="log10 + 0.1"; J="mg/L" },
    @{ Row=12; Name="mineral.N";  B=0.0150674108656426;    C=1.978743;   D="log10 + 0.1"; J="mg/L" },
    @{ Row=13; Name="Tot.N";      B=0.00838160366151803;   C=2.782443;   D="log10 + 0.1"; J="mg/L" },
    @{ Row=14; Name="PO4.P";      B=-0.00154712427761869;  C=0.06183478; D="log10 + 0.1"; J="mg/L" },
    @{ Row=15; Name="Tot.P";      B=-0.00238239917129794;  C=0.1049174;  D="log10 + 0.1"; J="mg/L" }
)

# Write column-by-column (matches the order new shared strings were
# originally interned: all Response names, then the two new SI units).
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 8).Value = 11
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 10).Value = $r.J
}

# ---- formulas --------------------------------------------------------
# E8, E9 keep their own one-off formula (mixed transform types among the
# new/blank rows), while E10:E15 share one formula block (all "log10+0.1").
$ws.Range("E8").Formula = "=B8/C8"
$ws.Range("E9").Formula = "=(10^B9-1)"
$ws.Range("E10:E15").Formula = "=(10^B10-1)"

# F3:F15 becomes one shared formula block covering the existing and new rows.
$ws.Range("F3:F15").Formula = "=E3*100"

# G8:G15 / I8:I15 are new shared formula blocks for the new rows.
$ws.Range("G8:G15").Formula = "=C8*(F8/100)"
$ws.Range("I8:I15").Formula = "=G8*H8"

Write-Host "formulas written"
